# Apply the "Chamada" (attendance) sheet edits:
#  - Add a new "Quarta" attendance column (G) for every student row,
#    mirroring the existing "Terça" column (F) which already held P/F marks.
#  - Normalize the style of the F column so it matches the plain
#    center/center style used elsewhere (drops a duplicate style def).
#  - Remove the now-unused trailing helper row (44); its helper cell
#    effectively becomes G43.
#  - Update the "active cell" selection bookmark.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chamada")

# Row -> attendance mark for the new "Quarta" column (G). Default is "P"
# (present); a handful of rows are marked "F" (absent), matching the
# pattern already present elsewhere on the same rows.
$gValues = @{}
foreach ($r in 4..22) { $gValues[$r] = "P" }
foreach ($r in 25..43) { $gValues[$r] = "P" }
$gValues[26] = "F"
$gValues[36] = "F"
$gValues[37] = "F"
$gValues[38] = "F"

# Rows 4-22 and 25-42 get a plain center-aligned G cell (same look as the
# existing F column). Row 43's G cell instead reuses the underlined
# "helper" style that used to live on F24 / L30 / old F44.
foreach ($r in 4..22) {
    $ws.Range("F" + $r).Font.Underline = $false
    $ws.Range("G" + $r).Value = $gValues[$r]
    $ws.Range("G" + $r).HorizontalAlignment = -4108
    $ws.Range("G" + $r).VerticalAlignment = -4108
}
foreach ($r in 25..42) {
    $ws.Range("F" + $r).Font.Underline = $false
    $ws.Range("G" + $r).Value = $gValues[$r]
    $ws.Range("G" + $r).HorizontalAlignment = -4108
    $ws.Range("G" + $r).VerticalAlignment = -4108
}

# Row 43: fix F43, then G43 takes the underlined helper style.
$ws.Range("F43").Font.Underline = $false
$ws.Range("G43").Value = $gValues[43]
$ws.Range("G43").HorizontalAlignment = -4108
$ws.Range("G43").VerticalAlignment = -4108
$ws.Range("G43").Font.Underline = 2

# Drop the old trailing helper row; its role is now played by G43.
$ws.Rows.Item(44).Delete()

# Move the remembered selection from the old helper cell to the new one.
$ws.Range("G43").Select()
